$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# For rows whose Price (column D) is being changed to a new value that
# looks like a plain number (e.g. "1.005", "26.365.57"), force the cell
# to Text format first so Excel stores the literal digits/dots instead
# of silently converting the string to a floating-point number.
# (Cells whose D value does not change are left completely untouched.)
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

$ws.Cells.Item(2, 2).Value2 = "Bitcoin"
$ws.Cells.Item(2, 3).Value2 = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).Value2 = "26.365.57"
$ws.Cells.Item(2, 5).Value2 = "  -1.94%  "

$ws.Cells.Item(3, 2).Value2 = "Ethereum"
$ws.Cells.Item(3, 3).Value2 = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).Value2 = "1.793.64"
$ws.Cells.Item(3, 5).Value2 = "  -1.78%  "

$ws.Cells.Item(4, 2).Value2 = "TetherUSD"
$ws.Cells.Item(4, 3).Value2 = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Cells.Item(4, 4).Value2 = "1.005"
$ws.Cells.Item(4, 5).Value2 = "  -0.20%  "

$ws.Cells.Item(5, 2).Value2 = "USDC"
$ws.Cells.Item(5, 3).Value2 = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(5, 4).Value2 = "1.005"
$ws.Cells.Item(5, 5).Value2 = "  -0.15%  "

$ws.Cells.Item(6, 2).Value2 = "BNB"
$ws.Cells.Item(6, 3).Value2 = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(6, 4).Value2 = "307.02"
$ws.Cells.Item(6, 5).Value2 = "  -1.24%  "

$ws.Cells.Item(7, 2).Value2 = "XRP"
$ws.Cells.Item(7, 3).Value2 = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(7, 4).Value2 = "0.4530"
$ws.Cells.Item(7, 5).Value2 = "  -1.10%  "

$ws.Cells.Item(8, 2).Value2 = "Cardano"
$ws.Cells.Item(8, 3).Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(8, 4).Value2 = "0.3591"
$ws.Cells.Item(8, 5).Value2 = "  -2.29%  "

$ws.Cells.Item(9, 2).Value2 = "OKB"
$ws.Cells.Item(9, 3).Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(9, 4).Value2 = "46.15"
$ws.Cells.Item(9, 5).Value2 = "  +0.45%  "

$ws.Cells.Item(10, 2).Value2 = "Dogecoin"
$ws.Cells.Item(10, 3).Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).Value2 = "0.07078"
$ws.Cells.Item(10, 5).Value2 = "  -1.05%  "

$ws.Cells.Item(11, 2).Value2 = "Polygon"
$ws.Cells.Item(11, 3).Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(11, 4).Value2 = "0.8862"
$ws.Cells.Item(11, 5).Value2 = "  +1.66%  "

$ws.Cells.Item(12, 2).Value2 = "TRON"
$ws.Cells.Item(12, 3).Value2 = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(12, 4).Value2 = "0.07766"
$ws.Cells.Item(12, 5).Value2 = "  -0.01%  "

$ws.Cells.Item(13, 2).Value2 = "Solana"
$ws.Cells.Item(13, 3).Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(13, 4).Value2 = "19.43"
$ws.Cells.Item(13, 5).Value2 = "  -0.52%  "

$ws.Cells.Item(14, 2).Value2 = "WrappedEther"
$ws.Cells.Item(14, 3).Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value2 = "1.815.69"
$ws.Cells.Item(14, 5).Value2 = "  -0.24%  "

$ws.Cells.Item(15, 2).Value2 = "Polkadot"
$ws.Cells.Item(15, 3).Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(15, 4).Value2 = "5.279"
$ws.Cells.Item(15, 5).Value2 = "  -0.64%  "

$ws.Cells.Item(16, 2).Value2 = "Chainlink"
$ws.Cells.Item(16, 3).Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(16, 4).Value2 = "6.321"
$ws.Cells.Item(16, 5).Value2 = "  -0.88%  "

$ws.Cells.Item(17, 2).Value2 = "Litecoin"
$ws.Cells.Item(17, 3).Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(17, 4).Value2 = "84.97"
$ws.Cells.Item(17, 5).Value2 = "  -2.19%  "

$ws.Cells.Item(18, 2).Value2 = "BinanceUSD"
$ws.Cells.Item(18, 3).Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(18, 4).Value2 = "1.007"
$ws.Cells.Item(18, 5).Value2 = "  -0.08%  "

$ws.Cells.Item(19, 2).Value2 = "ShibaInu"
$ws.Cells.Item(19, 3).Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).Value2 = "0.000008503"
$ws.Cells.Item(19, 5).Value2 = "  -2.30%  "

$ws.Cells.Item(20, 2).Value2 = "Dai"
$ws.Cells.Item(20, 3).Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(20, 4).Value2 = "1.005"
$ws.Cells.Item(20, 5).Value2 = "  -0.12%  "

$ws.Cells.Item(21, 2).Value2 = "WrappedBTC"
$ws.Cells.Item(21, 3).Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(21, 4).Value2 = "26.381.77"
$ws.Cells.Item(21, 5).Value2 = "  -1.98%  "

$ws.Cells.Item(22, 2).Value2 = "Avalanche"
$ws.Cells.Item(22, 3).Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(22, 4).Value2 = "14.21"
$ws.Cells.Item(22, 5).Value2 = "  -1.54%  "

$ws.Cells.Item(23, 2).Value2 = "Uniswap"
$ws.Cells.Item(23, 3).Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(23, 4).Value2 = "4.962"
$ws.Cells.Item(23, 5).Value2 = "  -0.36%  "

$ws.Cells.Item(24, 2).Value2 = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(24, 3).Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(24, 4).Value2 = "2.041.69"
$ws.Cells.Item(24, 5).Value2 = "  +0.08%  "

$ws.Cells.Item(25, 2).Value2 = "Cosmos"
$ws.Cells.Item(25, 3).Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(25, 4).Value2 = "10.56"
$ws.Cells.Item(25, 5).Value2 = "  +0.92%  "

$ws.Cells.Item(26, 2).Value2 = "Toncoin"
$ws.Cells.Item(26, 3).Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).Value2 = "1.968"
$ws.Cells.Item(26, 5).Value2 = "  -1.94%  "

$ws.Cells.Item(27, 2).Value2 = "Monero"
$ws.Cells.Item(27, 3).Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(27, 4).Value2 = "151.03"
$ws.Cells.Item(27, 5).Value2 = "  -0.07%  "

$ws.Cells.Item(28, 2).Value2 = "EthereumClassic"
$ws.Cells.Item(28, 3).Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).Value2 = "17.80"
$ws.Cells.Item(28, 5).Value2 = "  -2.16%  "

$ws.Cells.Item(29, 2).Value2 = "LidoDAOToken"
$ws.Cells.Item(29, 3).Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(29, 4).Value2 = "2.018"
$ws.Cells.Item(29, 5).Value2 = "  +2.85%  "

$ws.Cells.Item(30, 2).Value2 = "BitcoinCash"
$ws.Cells.Item(30, 3).Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(30, 4).Value2 = "111.80"
$ws.Cells.Item(30, 5).Value2 = "  -1.51%  "

$ws.Cells.Item(31, 2).Value2 = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).Value2 = "4.852"
$ws.Cells.Item(31, 5).Value2 = "  -1.38%  "

$ws.Cells.Item(32, 2).Value2 = "Stellar"
$ws.Cells.Item(32, 3).Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).Value2 = "0.08696"
$ws.Cells.Item(32, 5).Value2 = "  -1.14%  "

$ws.Cells.Item(33, 2).Value2 = "HuobiToken"
$ws.Cells.Item(33, 3).Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(33, 4).Value2 = "3.126"
$ws.Cells.Item(33, 5).Value2 = "  +1.48%  "

$ws.Cells.Item(34, 2).Value2 = "RenderToken"
$ws.Cells.Item(34, 3).Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(34, 4).Value2 = "2.783"
$ws.Cells.Item(34, 5).Value2 = "  +10.81%  "

$ws.Cells.Item(35, 2).Value2 = "Filecoin"
$ws.Cells.Item(35, 3).Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(35, 4).Value2 = "4.437"
$ws.Cells.Item(35, 5).Value2 = "  -0.80%  "

$ws.Cells.Item(36, 2).Value2 = "ImmutableX"
$ws.Cells.Item(36, 3).Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).Value2 = "0.7205"
$ws.Cells.Item(36, 5).Value2 = "  -3.31%  "

$ws.Cells.Item(37, 2).Value2 = "ARBITRUM"
$ws.Cells.Item(37, 3).Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).Value2 = "1.102"
$ws.Cells.Item(37, 5).Value2 = "  -2.19%  "

$ws.Cells.Item(38, 2).Value2 = "Frax"
$ws.Cells.Item(38, 3).Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(38, 4).Value2 = "1.004"
$ws.Cells.Item(38, 5).Value2 = "  +0.26%  "

$ws.Cells.Item(39, 2).Value2 = "TrustWalletToken"
$ws.Cells.Item(39, 3).Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(39, 4).Value2 = "1.067"
$ws.Cells.Item(39, 5).Value2 = "  -1.79%  "

$ws.Cells.Item(40, 2).Value2 = "VeChain"
$ws.Cells.Item(40, 3).Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).Value2 = "0.01930"
$ws.Cells.Item(40, 5).Value2 = "  -0.34%  "

$ws.Cells.Item(41, 2).Value2 = "Hedera"
$ws.Cells.Item(41, 3).Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(41, 4).Value2 = "0.05068"
$ws.Cells.Item(41, 5).Value2 = "  -0.87%  "

$ws.Cells.Item(42, 2).Value2 = "MXToken"
$ws.Cells.Item(42, 3).Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(42, 4).Value2 = "2.854"
$ws.Cells.Item(42, 5).Value2 = "  -1.54%  "

$ws.Cells.Item(43, 2).Value2 = "TheSandbox"
$ws.Cells.Item(43, 3).Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43, 4).Value2 = "0.5073"
$ws.Cells.Item(43, 5).Value2 = "  +2.26%  "

$ws.Cells.Item(44, 2).Value2 = "FraxShare"
$ws.Cells.Item(44, 3).Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).Value2 = "6.827"
$ws.Cells.Item(44, 5).Value2 = "  -1.39%  "

$ws.Cells.Item(45, 2).Value2 = "Algorand"
$ws.Cells.Item(45, 3).Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).Value2 = "0.1515"
$ws.Cells.Item(45, 5).Value2 = "  -5.13%  "

$ws.Cells.Item(46, 2).Value2 = "Aptos"
$ws.Cells.Item(46, 3).Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).Value2 = "8.021"
$ws.Cells.Item(46, 5).Value2 = "  -3.05%  "

$ws.Cells.Item(47, 2).Value2 = "PaxDollar"
$ws.Cells.Item(47, 3).Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(47, 4).Value2 = "1.005"
$ws.Cells.Item(47, 5).Value2 = "  -0.11%  "

$ws.Cells.Item(48, 2).Value2 = "Decentraland"
$ws.Cells.Item(48, 3).Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(48, 4).Value2 = "0.4633"
$ws.Cells.Item(48, 5).Value2 = "  -0.86%  "

$ws.Cells.Item(49, 2).Value2 = "EnergySwap"
$ws.Cells.Item(49, 3).Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value2 = "9.869"
$ws.Cells.Item(49, 5).Value2 = "  -2.36%  "

$ws.Cells.Item(50, 2).Value2 = "Quant"
$ws.Cells.Item(50, 3).Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(50, 4).Value2 = "100.81"
$ws.Cells.Item(50, 5).Value2 = "  -0.68%  "

$ws.Cells.Item(51, 2).Value2 = "NEARProtocol"
$ws.Cells.Item(51, 3).Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).Value2 = "1.564"
$ws.Cells.Item(51, 5).Value2 = "  -2.49%  "
